# New daily price-record row for "Perejil" (Vega Modelo de Temuco) is inserted
# at row 108, pushing the existing rows 108..173 down to 109..174.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(108).EntireRow.Insert()

$ws.Range("A108").Value = 10
$ws.Range("B108").Value = "Vega Modelo de Temuco"
$ws.Range("C108").Value = "La Araucanía"
$ws.Range("D108").Value = 44438
$ws.Range("E108").Value = 9
$ws.Range("F108").Value = 100112044
$ws.Range("G108").Value = "Perejil"
$ws.Range("H108").Value = "Sin especificar"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 40
$ws.Range("K108").Value = 4000
$ws.Range("L108").Value = 4000
$ws.Range("M108").Value = 4000
$ws.Range("N108").Value = "$/docena de atados (3 kilos)"
$ws.Range("O108").Value = "Provincia de Cautín"
$ws.Range("P108").Value = 1333
$ws.Range("Q108").Value = 3
$ws.Range("R108").Value = "Hortaliza"
